$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text (avoids Excel auto-converting
    # numeric-looking strings like "320.90" into the number 320.9), then
    # restore the cell's original (unstyled) formatting so we don't leave
    # a stray text-format style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "47.285.33"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.490.54"
$ws.Range("E3").Value = "  -0.11%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "320.90"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "108.30"
$ws.Range("E6").Value = "  +3.44%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.15%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "38.81"
$ws.Range("E10").Value = "  +5.37%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.22%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.02%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "18.33"
$ws.Range("E13").Value = "  +0.34%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.93%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.879.93"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.488.88"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +0.63%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "47.197.33"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D19") "12.95"
$ws.Range("E19").Value = "  +3.00%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +0.36%  "

# Row 21 - ShibaInu
Set-TextValue $ws.Range("D21") "0.0₃0933"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22 - ImmutableX
$ws.Range("E22").Value = "  +13.09%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "70.26"
$ws.Range("E23").Value = "  -0.36%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "245.08"
$ws.Range("E24").Value = "  -2.05%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.57"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "25.71"
$ws.Range("E27").Value = "  -1.50%  "

# Row 28 / Row 29 swap: Cosmos <-> Toncoin
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D28") "2.27"
$ws.Range("E28").Value = "  +3.65%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D29") "10.03"
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +3.32%  "

# Row 31 - InjectiveProtocol
Set-TextValue $ws.Range("D31") "34.62"
$ws.Range("E31").Value = "  -1.11%  "

# Row 32 - OKB
Set-TextValue $ws.Range("D32") "49.63"
$ws.Range("E32").Value = "  +0.39%  "

# Row 33 - Celestia
Set-TextValue $ws.Range("D33") "20.77"
$ws.Range("E33").Value = "  +6.34%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "5.35"
$ws.Range("E34").Value = "  +0.55%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0782"
$ws.Range("E35").Value = "  +1.09%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.11%  "

# Row 37 - RenderToken
Set-TextValue $ws.Range("D37") "4.75"
$ws.Range("E37").Value = "  +4.30%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +1.95%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws.Range("D39") "2.93"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40 - EnergySwap
Set-TextValue $ws.Range("D40") "23.33"
$ws.Range("E40").Value = "  +10.22%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -0.31%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +0.12%  "

# Row 43 - Monero
Set-TextValue $ws.Range("D43") "117.87"
$ws.Range("E43").Value = "  -3.27%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +0.30%  "

# Row 45 - Maker
Set-TextValue $ws.Range("D45") "1.990.39"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46 - NEARProtocol
$ws.Range("E46").Value = "  +2.07%  "

# Row 47 - ApeXProtocol
Set-TextValue $ws.Range("D47") "2.02"
$ws.Range("E47").Value = "  -5.42%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "9.12"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49 - Stacks
$ws.Range("E49").Value = "  -1.08%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  -5.45%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "56.66"
$ws.Range("E51").Value = "  +4.30%  "
